$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Increment the "Förändrad" (changed) date by one day (45174 -> 45175)
# for rows 2-5 in column C.
$ws.Range("C2").Value = 45175
$ws.Range("C3").Value = 45175
$ws.Range("C4").Value = 45175
$ws.Range("C5").Value = 45175
